$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 600 ("「言わぬが花」" entry) — all subsequent rows shift up by one.
$ws.Rows.Item(600).Delete()
